# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the crypto symbol
# list with the latest scraped values (GitHub Actions run of 2023-01-16).
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the sheet's existing inline-string cells) instead
# of auto-converting numeric-/percent-looking text into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'305.37"
$ws.Range("E2").Value = "'1.61%"
$ws.Range("D3").Value = "'32.20"
$ws.Range("E3").Value = "'2.05%"
$ws.Range("D4").Value = "'5.220"
$ws.Range("E4").Value = "'1.54%"
$ws.Range("D5").Value = "'0.07479"
$ws.Range("E5").Value = "'-0.12%"
$ws.Range("D6").Value = "'2.371"
$ws.Range("E6").Value = "'43.58%"
$ws.Range("D7").Value = "'7.996"
$ws.Range("E7").Value = "'1.96%"
$ws.Range("D8").Value = "'3.865"
$ws.Range("D9").Value = "'0.9174"
$ws.Range("E9").Value = "'-0.70%"
$ws.Range("D10").Value = "'0.1734"
$ws.Range("E10").Value = "'1.06%"
$ws.Range("D11").Value = "'0.07719"
$ws.Range("E11").Value = "'1.04%"
$ws.Range("D12").Value = "'0.08246"
$ws.Range("E12").Value = "'2.79%"
$ws.Range("D13").Value = "'0.03016"
$ws.Range("E13").Value = "'0.49%"
$ws.Range("D14").Value = "'0.09947"
$ws.Range("E14").Value = "'0.55%"
$ws.Range("D15").Value = "'0.001503"
$ws.Range("E15").Value = "'0.25%"
$ws.Range("D16").Value = "'0.006137"
$ws.Range("E16").Value = "'-1.18%"
$ws.Range("D17").Value = "'3.493"
$ws.Range("E17").Value = "'1.40%"
$ws.Range("E18").Value = "'-0.08%"
$ws.Range("D19").Value = "'0.3265"
$ws.Range("E19").Value = "'-0.86%"
$ws.Range("D20").Value = "'0.1346"
$ws.Range("E20").Value = "'0.64%"
$ws.Range("D21").Value = "'4.646"
$ws.Range("E21").Value = "'1.35%"
$ws.Range("E22").Value = "'-1.29%"
$ws.Range("D23").Value = "'0.1561"
$ws.Range("E23").Value = "'0.58%"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'1.46%"
$ws.Range("D25").Value = "'0.004531"
$ws.Range("E25").Value = "'2.34%"
$ws.Range("D26").Value = "'0.0001295"
$ws.Range("E26").Value = "'-7.58%"
$ws.Range("D27").Value = "'0.0002732"
$ws.Range("E27").Value = "'51.95%"
$ws.Range("D39").Value = "'0.01784"
$ws.Range("E39").Value = "'7.87%"
$ws.Range("D40").Value = "'0.04575"
$ws.Range("E40").Value = "'0.68%"
$ws.Range("D41").Value = "'0.007374"
$ws.Range("E41").Value = "'6.07%"
$ws.Range("D42").Value = "'0.1363"
$ws.Range("E42").Value = "'1.50%"
$ws.Range("D43").Value = "'0.002172"
$ws.Range("E43").Value = "'5.33%"
$ws.Range("E44").Value = "'-13.06%"
$ws.Range("D45").Value = "'0.00006465"
$ws.Range("E45").Value = "'6.33%"
$ws.Range("E46").Value = "'15.26%"
$ws.Range("D47").Value = "'0.009864"
$ws.Range("E47").Value = "'-19.51%"
